$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

# --- Burndown data update (Day 3 / column H) ---------------------------
# The "Skattat" (Estimated) row (row 2) and the chart's source summary row
# (row 6) had no remaining/ongoing work recorded for Day 3; the author
# worked on the project today and updated the burndown numbers to 5.
$ws.Range("H2").Value = 5
$ws.Range("H6").Value = 5

# --- Keep the burndown chart's "Skattat" series pointing at the same ----
# source range so it picks up the refreshed Day 3 figure above.
$co = $ws.ChartObjects(1)
$chart = $co.Chart
$series = $chart.SeriesCollection(1)
$series.Formula = "=SERIES(Blad1!`$E`$6,,Blad1!`$F`$6:`$J`$6,1)"

# --- Restore the author's last selection on the sheet -------------------
$ws.Range("I9").Select()
